$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'47.595.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.05%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.493.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'322.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.14%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'109.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +3.62%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.523"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.56%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.03%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.05%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'39.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +4.16%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.46%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +0.80%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +2.19%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +0.69%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.882.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.16%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.492.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.07%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.849"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.75%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'47.443.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.03%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'13.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +5.95%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +1.21%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.57%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +15.47%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'70.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.02%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'247.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.48%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.63%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.11%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'25.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.39%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("B28").Value = "'Toncoin"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'2.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +4.53%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "'Cosmos"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'10.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.12%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.140"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +3.68%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'34.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.27%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +0.96%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'20.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +2.76%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'5.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.36%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +0.90%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +0.15%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +2.61%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +1.15%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.56%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.36%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'22.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +6.48%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'2.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.77%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'119.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.77%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +0.00%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.997.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.76%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +2.33%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -2.56%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.08%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'9.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.33%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'5.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.39%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'56.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +3.43%  "
$ws.Range("E51").Style = "Normal"
